$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update frequency counts for unchanged category pairs (rows 2-7)
$ws.Range("C2").Value = 3876
$ws.Range("C3").Value = 3704
$ws.Range("C4").Value = 2887
$ws.Range("C5").Value = 1940
$ws.Range("C6").Value = 1749
$ws.Range("C7").Value = 806

# Row 8: category pair & frequency change
$ws.Range("A8").Value = "Textiles & Cozy Items"
$ws.Range("B8").Value = "Textiles & Cozy Items"
$ws.Range("C8").Value = 589

# Row 9: category pair & frequency change
$ws.Range("A9").Value = "Kitchen & Dining"
$ws.Range("B9").Value = "Home Decor"
$ws.Range("C9").Value = 556

# Row 10: complementary category & frequency change
$ws.Range("B10").Value = "Seasonal & Holidays"
$ws.Range("C10").Value = 508

# Row 11: complementary category & frequency change
$ws.Range("B11").Value = "Vintage & Collectibles"
$ws.Range("C11").Value = 495
